$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C30").Value = "Volatge Regulator Options"
$ws.Range("D30").Value = "Mouser lInk"
$ws.Range("E30").Value = "Digi Link"

$ws.Range("C31").Value = "TPS78533QWDRBRQ1"
$ws.Range("D31").Value = "https://www.mouser.ca/ProductDetail/Texas-Instruments/TPS78533QWDRBRQ1?qs=iLbezkQI%252BsiNbdqh2Ko9GQ%3D%3D"
$ws.Range("E31").Value = "https://www.digikey.ca/en/products/detail/texas-instruments/tps78533qwdrbrq1/14123967"

$ws.Range("C32").Value = "LM1117MP-3.3/NOPB"
$ws.Range("D32").Value = "https://www.mouser.ca/ProductDetail/Texas-Instruments/LM1117MP-3.3-NOPB?qs=X1J7HmVL2ZFn4x9DZ4T2hA%3D%3D"
$ws.Range("E32").Value = "https://www.digikey.ca/en/products/detail/texas-instruments/LM1117MP-3-3-NOPB/304882"

$ws.Range("C33").Value = "LDL1117S33R"
$ws.Range("D33").Value = "https://www.mouser.ca/ProductDetail/STMicroelectronics/LDL1117S33R?qs=AQlKX63v8Rt9Bf6AWSrbFg%3D%3D"
$ws.Range("E33").Value = "https://www.digikey.ca/en/products/detail/stmicroelectronics/LDL1117S33R/7102071"

$ws.Range("C31").Font.Name = "Arial"
$ws.Range("C31").Font.Size = 13
$ws.Range("C31").Font.Color = 3355443

$ws.Range("G39").Select() | Out-Null
